$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.512.39'

$ws.Range('D3').Value = '1.619.03'
$ws.Range('E3').Value = '  -1.65%  '

$ws.Range('E4').Value = '  +0.04%  '

$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '210.92'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -1.09%  '

$ws.Range('E6').Value = '  -1.68%  '

$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '22.84'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  -1.22%  '

$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.261'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  +1.15%  '

$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.0613'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  -0.08%  '

$ws.Range('E11').Value = '  -0.51%  '

$ws.Range('D12').Value = '1.848.45'
$ws.Range('E12').Value = '  -1.63%  '

$ws.Range('D13').Value = '1.602.39'
$ws.Range('E13').Value = '  -2.61%  '

$ws.Range('E14').Value = '  -0.22%  '

$ws.Range('E15').Value = '  -2.06%  '

$ws.Range('E16').Value = '  +1.40%  '

$ws.Range('D17').Value = '27.492.46'
$ws.Range('E17').Value = '  -0.74%  '

$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '230.84'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  -0.08%  '

$ws.Range('D19').Value = '0.0₃0718'
$ws.Range('E19').Value = '  -1.04%  '

$ws.Range('E20').Value = '  -2.06%  '

$ws.Range('E22').Value = '  -0.86%  '

$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '10.17'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +1.48%  '

$ws.Range('E24').Value = '  +5.55%  '

$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '150.89'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +0.96%  '

$ws.Range('E27').Value = '  -1.98%  '

$ws.Range('E28').Value = '  +0.04%  '

$ws.Range('E29').Value = '  -0.96%  '

$ws.Range('E30').Value = '  -0.89%  '

$ws.Range('E31').Value = '  -0.73%  '

$ws.Range('E32').Value = '  -1.21%  '

$ws.Range('D33').Value = '1.468.03'
$ws.Range('E33').Value = '  +1.81%  '

$ws.Range('E34').Value = '  -3.02%  '

$ws.Range('E35').Value = '  -3.68%  '

$ws.Range('E36').Value = '  -0.10%  '

$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '0.950'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  +5.46%  '

$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.0167'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  -0.28%  '

$ws.Range('B39').Value = 'ImmutableX'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.559'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -2.02%  '

$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '0.859'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  -2.88%  '

$ws.Range('E41').Value = '  +0.04%  '

$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '67.82'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  +3.35%  '

$ws.Range('B43').Value = 'WEMIXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '0.988'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  -4.25%  '

$ws.Range('B44').Value = 'MXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '2.21'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  -2.37%  '

$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '5.25'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -7.24%  '

$ws.Range('B46').Value = 'RocketPoolETH'
$ws.Range('C46').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D46').Value = '1.759.43'
$ws.Range('E46').Value = '  -1.57%  '

$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '1.71'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +0.97%  '

$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '86.59'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +0.76%  '

$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.0₆0105'
$ws.Range('E49').Value = '  -2.43%  '

$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '0.101'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  +1.92%  '

$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '7.65'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -1.42%  '
